$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data series shifts forward by one decade: rows that held 2000-2009
# (old rows 2-11) are dropped, rows holding 2010-2020 (old rows 12-22)
# become the new rows 2-12, and two new rows are appended for 2021/2022.
# Net effect: delete the surplus trailing rows 15-22, then rewrite rows
# 2-14 with the new year labels and values.

$ws.Rows("15:22").Delete() | Out-Null

$data = @(
    @("2010年", 37.9, 31.9, 33.4),
    @("2011年", 37.1, 32.3, 33.6),
    @("2012年", 35.9, 32, 33),
    @("2013年", 34.1, 30.1, 31.2),
    @("2014年", 33.6, 30, 31),
    @("2015年", 33, 29.7, 30.6),
    @("2016年", 32.2, 29.3, 30.1),
    @("2017年", 31.2, 28.6, 29.3),
    @("2018年", 30.1, 27.7, 28.4),
    @("2019年", 30, 27.6, 28.2),
    @("2020年", 32.7, 29.2, 30.2),
    @("2021年", 32.7, 28.6, 29.8),
    @("2022年", 32.9806665375964, 29.4771654602819, 30.4872906929388)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
